# Chiffres COVID-19 Valais - daily data update
# Updates "Nb nouveaux cas positifs" (C), intubated/SI/hospit counts (E,F,G)
# and new-death-location splits (L,M) for a handful of existing rows, and
# fills in three previously-blank rows (319-321) with real data. All the
# other changed cells in the target file (B, H, J, K and the cumulative
# running totals) are formulas that recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 288 : nouveaux cas 144 -> 147 ---
$ws.Range("C288").Value2 = 147

# --- Row 315 : nouveaux cas 170 -> 169 ---
$ws.Range("C315").Value2 = 169

# --- Row 316 : nouveaux deces extra-hospitaliers 2 -> 3 ---
$ws.Range("M316").Value2 = 3

# --- Row 317 : nouveaux cas 115 -> 145, deces hopital 0 -> 1, extra-hosp 0 -> 3 ---
$ws.Range("C317").Value2 = 145
$ws.Range("L317").Value2 = 1
$ws.Range("M317").Value2 = 3

# --- Row 318 : nouveaux cas 21 -> 147, deces hopital 0 -> 3, extra-hosp 0 -> 1 ---
$ws.Range("C318").Value2 = 147
$ws.Range("L318").Value2 = 3
$ws.Range("M318").Value2 = 1

# --- Row 319 (2021-01-09) : first entry of previously-empty row ---
$ws.Range("C319").Value2 = 68
$ws.Range("E319").Value2 = 12
$ws.Range("F319").Value2 = 7
$ws.Range("G319").Value2 = 83
$ws.Range("L319").Value2 = 3
$ws.Range("M319").Value2 = 1

# --- Row 320 (2021-01-10) ---
$ws.Range("C320").Value2 = 46
$ws.Range("E320").Value2 = 12
$ws.Range("F320").Value2 = 7
$ws.Range("G320").Value2 = 87
$ws.Range("L320").Value2 = 1
$ws.Range("M320").Value2 = 1

# --- Row 321 (2021-01-11) ---
$ws.Range("C321").Value2 = 26
$ws.Range("E321").Value2 = 10
$ws.Range("F321").Value2 = 7
$ws.Range("G321").Value2 = 87
$ws.Range("L321").Value2 = 0
$ws.Range("M321").Value2 = 0

# --- Selection cosmetic change: active cell moves back to A2 ---
$ws.Activate()
$ws.Range("A2").Select()
